# Regenerate save_data column G ("K") values.
# The commit replaces the old "Strike#" derived values in column G with
# newly (re)calculated "K" values for each data row (rows 2-34 on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    8  = 2
    9  = 2
    10 = 2
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 3
    17 = 3
    18 = 3
    19 = 0
    20 = 0
    21 = 1
    22 = 3
    23 = 2
    24 = 2
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 2
    30 = 3
    31 = 0
    32 = 1
    33 = 0
    34 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
